# Update the division problems in the single table of the worksheet.
# Only the five populated rows (1, 5, 9, 13, 17) hold content; the other
# rows are blank "answer" rows. Cells are addressed by (row, column) so
# that duplicate values (e.g. "32÷9=" appears twice in the source) are
# replaced independently with their correct targets.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; Text = "55÷5=" },
    @{ Row = 1;  Col = 2; Text = "99÷9=" },
    @{ Row = 1;  Col = 3; Text = "17÷2=" },
    @{ Row = 1;  Col = 4; Text = "28÷7=" },
    @{ Row = 1;  Col = 5; Text = "67÷6=" },

    @{ Row = 5;  Col = 1; Text = "24÷9=" },
    @{ Row = 5;  Col = 2; Text = "81÷5=" },
    @{ Row = 5;  Col = 3; Text = "94÷3=" },
    @{ Row = 5;  Col = 4; Text = "67÷4=" },
    @{ Row = 5;  Col = 5; Text = "12÷7=" },

    @{ Row = 9;  Col = 1; Text = "65÷4=" },
    @{ Row = 9;  Col = 2; Text = "62÷9=" },
    @{ Row = 9;  Col = 3; Text = "36÷6=" },
    @{ Row = 9;  Col = 4; Text = "39÷3=" },
    @{ Row = 9;  Col = 5; Text = "24÷2=" },

    @{ Row = 13; Col = 1; Text = "87÷9=" },
    @{ Row = 13; Col = 2; Text = "74÷5=" },
    @{ Row = 13; Col = 3; Text = "26÷9=" },
    @{ Row = 13; Col = 4; Text = "32÷9=" },
    @{ Row = 13; Col = 5; Text = "51÷8=" },

    @{ Row = 17; Col = 1; Text = "22÷8=" },
    @{ Row = 17; Col = 2; Text = "80÷5=" },
    @{ Row = 17; Col = 3; Text = "31÷2=" },
    @{ Row = 17; Col = 4; Text = "45÷3=" },
    @{ Row = 17; Col = 5; Text = "28÷5=" }
)

foreach ($u in $updates) {
    $cell = $t.Rows.Item($u.Row).Cells.Item($u.Col)
    $cell.Range.Text = $u.Text
}
